$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows above row 123, pushing the existing 123:128 block
# (the previous week's data, dated 44509) down to become rows 129:134.
$ws.Rows("123:128").Insert()

# Common values shared by every row in this block.
$mercadoId = 9
$mercado   = "Vega Central Mapocho de Santiago"
$region    = "Metropolitana"
$fecha     = 44516
$codreg    = 13
$tipo      = "Fruta"
$productoId = 100107
$producto  = "Otros"
$categoriaId = 100107002
$categoria = "Chirimoya"
$variedad  = "Cultivar IV Región"
$origen    = "Provincia de Limarí"

# Per-row (quality-specific) data for the new week (row 123..128).
$data = @(
    @{ Row=123; Calidad="Cuarta";                  Volumen=250; PMin=1000;  PMax=1000;  PProm=1000;  Unidad="$/kilo (en caja de 15 kilos)"; PrecioKg=1000; KgUnidad=1 },
    @{ Row=124; Calidad="Especial";                Volumen=220; PMin=20000; PMax=20000; PProm=20000; Unidad="$/bandeja 8 kilos";             PrecioKg=2500; KgUnidad=8 },
    @{ Row=125; Calidad="Extra (doble especial)";  Volumen=310; PMin=24000; PMax=24000; PProm=24000; Unidad="$/bandeja 8 kilos";             PrecioKg=3000; KgUnidad=8 },
    @{ Row=126; Calidad="Primera";                  Volumen=280; PMin=16000; PMax=16000; PProm=16000; Unidad="$/bandeja 8 kilos";             PrecioKg=2000; KgUnidad=8 },
    @{ Row=127; Calidad="Segunda";                  Volumen=300; PMin=12800; PMax=12800; PProm=12800; Unidad="$/bandeja 8 kilos";             PrecioKg=1600; KgUnidad=8 },
    @{ Row=128; Calidad="Tercera";                  Volumen=220; PMin=1400;  PMax=1400;  PProm=1400;  Unidad="$/kilo (en caja de 15 kilos)"; PrecioKg=1400; KgUnidad=1 }
)

foreach ($d in $data) {
    $r = $d.Row
    $ws.Cells.Item($r, 1).Value2 = $mercadoId
    $ws.Cells.Item($r, 2).Value2 = $mercado
    $ws.Cells.Item($r, 3).Value2 = $region
    $ws.Cells.Item($r, 4).Value2 = $fecha
    $ws.Cells.Item($r, 5).Value2 = $codreg
    $ws.Cells.Item($r, 6).Value2 = $tipo
    $ws.Cells.Item($r, 7).Value2 = $productoId
    $ws.Cells.Item($r, 8).Value2 = $producto
    $ws.Cells.Item($r, 9).Value2 = $categoriaId
    $ws.Cells.Item($r, 10).Value2 = $categoria
    $ws.Cells.Item($r, 11).Value2 = $variedad
    $ws.Cells.Item($r, 12).Value2 = $d.Calidad
    $ws.Cells.Item($r, 13).Value2 = $d.Volumen
    $ws.Cells.Item($r, 14).Value2 = $d.PMin
    $ws.Cells.Item($r, 15).Value2 = $d.PMax
    $ws.Cells.Item($r, 16).Value2 = $d.PProm
    $ws.Cells.Item($r, 17).Value2 = $d.Unidad
    $ws.Cells.Item($r, 18).Value2 = $origen
    $ws.Cells.Item($r, 19).Value2 = $d.PrecioKg
    $ws.Cells.Item($r, 20).Value2 = $d.KgUnidad
}
